$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 through 10 (the old Table[2] and Table[3] rows)
$ws.Rows("5:10").Delete()

# Update remaining data rows to reflect the new single-table XPath set
$ws.Range("B2").Value = "/NewDataSet/Table[1]/Town"
$ws.Range("C2").Value = "[A-Z a-z].*"

$ws.Range("B3").Value = "/NewDataSet/Table[1]/County"
$ws.Range("C3").Value = "[A-Z a-z 0-9].*"

$ws.Range("B4").Value = "/NewDataSet/Table[1]/PostCode"
$ws.Range("C4").Value = "[A-Z a-z 0-9].*"

# Move active selection to B5 as in the target file
$ws.Range("B5").Select()
